$wb = $excel.ActiveWorkbook

# 1. Rename second sheet
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsInclude = $wb.Worksheets.Item("Include from Ferlab.bio CodeS")
$wsInclude.Name = "Include #0"

# 2. Update Metadata sheet cells
# Date value changes
$wsMeta.Range("B8").Value = "2024-10-02T15:04:17+00:00"
# Contact value changes
$wsMeta.Range("B10").Value = "Ferlab.bio (http://example.org/example-publisher)"

# 3. Insert new row "Jurisdiction" after row 10 (Contact), before Description
$wsMeta.Range("A11:B11").Insert()
$wsMeta.Range("A11").Value = "Jurisdiction"
$wsMeta.Range("B11").Value = ""

# Copy formatting from row above (Contact row) for the new row cells
$wsMeta.Range("A10").Copy()
$wsMeta.Range("A11").PasteSpecial(-4122) | Out-Null
$wsMeta.Range("B10").Copy()
$wsMeta.Range("B11").PasteSpecial(-4122) | Out-Null

# 4. Add new "Immutable" row at the bottom (row 15) with value "BooleanType[null]"
$wsMeta.Range("A15").Value = "Immutable"
$wsMeta.Range("B15").Value = "BooleanType[null]"

$wsMeta.Range("A14").Copy()
$wsMeta.Range("A15").PasteSpecial(-4122) | Out-Null
$wsMeta.Range("B14").Copy()
$wsMeta.Range("B15").PasteSpecial(-4122) | Out-Null
